# Update cell values in the "Inscricoes" worksheet (Table1) per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    5 = @{ "F" = 88; "H" = 88 }
    10 = @{ "E" = 472; "F" = 233; "H" = 233 }
    11 = @{ "F" = 178; "H" = 178 }
    12 = @{ "E" = 461; "F" = 253; "H" = 253 }
    14 = @{ "F" = 62; "H" = 62 }
    15 = @{ "E" = 151 }
    16 = @{ "E" = 187 }
    18 = @{ "E" = 50 }
    21 = @{ "E" = 132; "F" = 72; "H" = 72 }
    22 = @{ "F" = 83; "H" = 83 }
    23 = @{ "E" = 190; "F" = 83; "H" = 83 }
    24 = @{ "E" = 198; "F" = 107; "H" = 107 }
    25 = @{ "E" = 248; "F" = 118; "H" = 118 }
    26 = @{ "E" = 149; "F" = 90; "H" = 90 }
    27 = @{ "E" = 307; "F" = 150; "H" = 150 }
    28 = @{ "F" = 70; "H" = 70 }
    29 = @{ "F" = 89; "H" = 89 }
    30 = @{ "F" = 116; "H" = 116 }
    31 = @{ "F" = 31; "H" = 31 }
    32 = @{ "F" = 103; "H" = 103 }
    33 = @{ "E" = 270; "F" = 138; "H" = 138 }
    34 = @{ "E" = 202; "F" = 130; "H" = 130 }
    35 = @{ "E" = 139 }
    37 = @{ "E" = 149; "F" = 74; "H" = 74 }
    39 = @{ "E" = 171; "F" = 81; "H" = 81 }
    40 = @{ "F" = 112; "H" = 112 }
    41 = @{ "E" = 377; "F" = 173; "H" = 173 }
    42 = @{ "E" = 344; "F" = 186; "H" = 186 }
    44 = @{ "F" = 147; "H" = 147 }
    45 = @{ "E" = 132; "F" = 67; "H" = 67 }
    46 = @{ "E" = 294; "F" = 161; "H" = 161 }
    47 = @{ "F" = 205; "H" = 205 }
    48 = @{ "E" = 191; "F" = 79; "H" = 79 }
    50 = @{ "E" = 231; "F" = 106; "H" = 106 }
    51 = @{ "F" = 92; "H" = 92 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

$wb.Save()
